$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.770.20'
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').Value = '1.613.82'
$ws.Range('E3').Value = '  -3.76%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '208.32'
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5183'
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.006'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2551'
$ws.Range('E8').Value = '  -4.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06177'
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.02'
$ws.Range('E10').Value = '  -6.54%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07526'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').Value = '1.620.18'
$ws.Range('E12').Value = '  -3.41%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.328'
$ws.Range('E13').Value = '  -3.17%  '
$ws.Range('D14').Value = '1.844.62'
$ws.Range('E14').Value = '  -3.19%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5387'
$ws.Range('E15').Value = '  -4.35%  '
$ws.Range('D16').Value = '0.0₅7818'
$ws.Range('E16').Value = '  -2.86%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.98'
$ws.Range('E17').Value = '  -4.63%  '
$ws.Range('D18').Value = '25.764.70'
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.573'
$ws.Range('E20').Value = '  -5.26%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '182.99'
$ws.Range('E21').Value = '  -2.90%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.962'
$ws.Range('E22').Value = '  -4.46%  '
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.002'
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.55'
$ws.Range('E25').Value = '  -3.53%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1200'
$ws.Range('E26').Value = '  -4.33%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.276'
$ws.Range('E27').Value = '  -4.09%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.39'
$ws.Range('E28').Value = '  -4.24%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05872'
$ws.Range('E30').Value = '  -5.34%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.236'
$ws.Range('E31').Value = '  -3.90%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.321'
$ws.Range('E32').Value = '  -5.24%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.304'
$ws.Range('E33').Value = '  -4.11%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.579'
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9581'
$ws.Range('E35').Value = '  -4.46%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.394'
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.707'
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5635'
$ws.Range('E38').Value = '  -7.21%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01576'
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.005'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8320'
$ws.Range('E41').Value = '  -4.48%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.012.79'
$ws.Range('E42').Value = '  -6.63%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.562'
$ws.Range('E43').Value = '  -8.76%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '98.77'
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('D45').Value = '1.769.05'
$ws.Range('E45').Value = '  -3.11%  '
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.009'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '53.77'
$ws.Range('E48').Value = '  -4.19%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05159'
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.829'
$ws.Range('E50').Value = '  -2.50%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4217'
$ws.Range('E51').Value = '  -0.76%  '
